# Update countries & provincias Spain
# - Refresh the "Datos actualizados..." timestamp (11:03 -> 12:03)
# - Swap the Montserrat / Seychelles rows (source list re-sorted) and
#   refresh their stats
# - Refresh the daily COVID stats for several country rows and the two
#   Spain "provincias" rows near the bottom of the sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 6 de Mayo de 2020 a las 12:03"

# --- Country rows (B:Casos totales, C:Nuevos casos, D:Casos activos,
#     E:Recuperados, F:Casos criticos, G:Muertes hoy, H:Muertes) -------

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 1238040
$ws.Range("C4").Value = 407
$ws.Range("E4").Value = 965087
$ws.Range("G4").Value = 13
$ws.Range("H4").Value = 72284

# Row 17 - Belgica
$ws.Range("B17").Value = 50781
$ws.Range("C17").Value = 272
$ws.Range("D17").Value = 12731
$ws.Range("E17").Value = 29711
$ws.Range("G17").Value = 323
$ws.Range("H17").Value = 8339

# Row 30 - Bielorrusia
$ws.Range("B30").Value = 19255
$ws.Range("C30").Value = 905
$ws.Range("D30").Value = 4388
$ws.Range("E30").Value = 14755
$ws.Range("G30").Value = 5
$ws.Range("H30").Value = 112

# Row 37 - Rumania
$ws.Range("E37").Value = 7525
$ws.Range("G37").Value = 17
$ws.Range("H37").Value = 858

# Row 39 - Indonesia
$ws.Range("B39").Value = 12438
$ws.Range("C39").Value = 367
$ws.Range("D39").Value = 2317
$ws.Range("E39").Value = 9226
$ws.Range("G39").Value = 23
$ws.Range("H39").Value = 895

# Row 40 - Banglades
$ws.Range("B40").Value = 11719
$ws.Range("C40").Value = 790
$ws.Range("E40").Value = 10130
$ws.Range("G40").Value = 3
$ws.Range("H40").Value = 186

# Row 48 - Serbia
$ws.Range("B48").Value = 7899
$ws.Range("C48").Value = 3
$ws.Range("D48").Value = 4017
$ws.Range("E48").Value = 3624
$ws.Range("G48").Value = 1
$ws.Range("H48").Value = 258

# Row 55 - Chequia
$ws.Range("B55").Value = 5573
$ws.Range("C55").Value = 161
$ws.Range("E55").Value = 1827

# --- Montserrat / Seychelles: swap order + refresh their own stats ----
# Row 205 used to be Montserrat, now becomes Seychelles
$ws.Range("A205").Value = "Seychelles"
$ws.Range("D205").Value = 8
$ws.Range("F205").Value = 0
$ws.Range("H205").Value = 0

# Row 206 used to be Seychelles, now becomes Montserrat
$ws.Range("A206").Value = "Montserrat"
$ws.Range("D206").Value = 7
$ws.Range("F206").Value = 1
$ws.Range("H206").Value = 1
